# Generate Report for Handoff
#
# Two new localization source files were picked up by the pipeline:
#   42f4766d-43ab-4701-a29c-f352d6e86e17.md
#   b5df1cc1-dfc1-4509-9a88-1b812e811937.md
# Both are "Ready for handoff" and need to be reported just above the
# pre-existing f0ea69a4-... row (which keeps all of its original data,
# just shifted down two rows) on every worksheet: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# Helper color/underline constants matching the workbook's existing
# "HyperLink" cell style (font color FF6495ED, single underline).
$hyperlinkColor = 15570276   # OLE (BGR) form of RGB 6495ED
$hyperlinkUnderline = 2      # xlUnderlineStyleSingle

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Make room: push the existing last row (row 3, f0ea69a4...) down by two
# rows, carrying its formatting with it.
$ws1.Rows.Item(3).Insert()
$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "42f4766d-43ab-4701-a29c-f352d6e86e17.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-41-13 10:41:01"

$ws1.Range("A4").Value = "b5df1cc1-dfc1-4509-9a88-1b812e811937.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-41-13 10:41:01"

# Rebuild the File Name hyperlinks top to bottom so the relationship ids
# line up in row order (A2, A3, A4, A5).
$ws1.UsedRange.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b7338fff35708235e71d882e94de97c35b7ecac/e2e/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md", "", "", "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/42f4766d43ab4701a29cf352d6e86e170000000/e2e/42f4766d-43ab-4701-a29c-f352d6e86e17.md", "", "", "42f4766d-43ab-4701-a29c-f352d6e86e17.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b5df1cc1dfc145099a881b812e8119370000000/e2e/b5df1cc1-dfc1-4509-9a88-1b812e811937.md", "", "", "b5df1cc1-dfc1-4509-9a88-1b812e811937.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/8b134427431c13d284d55392df0aa999e8a1ea98/e2e/f0ea69a4-53a8-458c-bc21-83bb0ca46deb.md", "", "", "f0ea69a4-53a8-458c-bc21-83bb0ca46deb.md") | Out-Null

# Keep the File Name column's look consistent with the workbook's
# existing hyperlink styling.
$ws1.Range("A2:A5").Font.Underline = $hyperlinkUnderline
$ws1.Range("A2:A5").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(3).Insert()
$ws2.Rows.Item(3).Insert()

$ws2.Range("A3").Value = "42f4766d-43ab-4701-a29c-f352d6e86e17.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "42f4766d-43ab-4701-a29c-f352d6e86e17.83e06547b6f1b7686a5c531cdcbf212a906eda13.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-13 10:40:57"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"

$ws2.Range("A4").Value = "b5df1cc1-dfc1-4509-9a88-1b812e811937.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "b5df1cc1-dfc1-4509-9a88-1b812e811937.e97104f063c9fdfd4fc8a5077fc084d94d44698e.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-13 10:40:57"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("I4").Value = "Include"

$ws2.UsedRange.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b7338fff35708235e71d882e94de97c35b7ecac/e2e/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md", "", "", "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b7338fff35708235e71d882e94de97c35b7ecac/e2e/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/61068100a525372f0aeb6d0e1c3fd8988dd6526d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.zh-cn.xlf", "", "", "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/58ae7ae447053e002200596ccb8c17cafbf5e87e/e2e/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md", "", "", "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4767e72dcd7524ffb76215bee09081b43d70d9db/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.zh-cn.xlf", "", "", "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/42f4766d43ab4701a29cf352d6e86e170000000/e2e/42f4766d-43ab-4701-a29c-f352d6e86e17.md", "", "", "42f4766d-43ab-4701-a29c-f352d6e86e17.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/42f4766d43ab4701a29cf352d6e86e170000000/e2e/42f4766d-43ab-4701-a29c-f352d6e86e17.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83e06547b6f1b7686a5c531cdcbf212a906eda13/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/42f4766d-43ab-4701-a29c-f352d6e86e17.83e06547b6f1b7686a5c531cdcbf212a906eda13.zh-cn.xlf", "", "", "42f4766d-43ab-4701-a29c-f352d6e86e17.83e06547b6f1b7686a5c531cdcbf212a906eda13.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b5df1cc1dfc145099a881b812e8119370000000/e2e/b5df1cc1-dfc1-4509-9a88-1b812e811937.md", "", "", "b5df1cc1-dfc1-4509-9a88-1b812e811937.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/b5df1cc1dfc145099a881b812e8119370000000/e2e/b5df1cc1-dfc1-4509-9a88-1b812e811937.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e97104f063c9fdfd4fc8a5077fc084d94d44698e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b5df1cc1-dfc1-4509-9a88-1b812e811937.e97104f063c9fdfd4fc8a5077fc084d94d44698e.zh-cn.xlf", "", "", "b5df1cc1-dfc1-4509-9a88-1b812e811937.e97104f063c9fdfd4fc8a5077fc084d94d44698e.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/8b134427431c13d284d55392df0aa999e8a1ea98/e2e/f0ea69a4-53a8-458c-bc21-83bb0ca46deb.md", "", "", "f0ea69a4-53a8-458c-bc21-83bb0ca46deb.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/8b134427431c13d284d55392df0aa999e8a1ea98/e2e/f0ea69a4-53a8-458c-bc21-83bb0ca46deb.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b536ec2138909cb5dcde27057c2767ec49b88d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f0ea69a4-53a8-458c-bc21-83bb0ca46deb.91b5d2b6863b784fabb2a985a9195c04a47ab4b9.zh-cn.xlf", "", "", "f0ea69a4-53a8-458c-bc21-83bb0ca46deb.91b5d2b6863b784fabb2a985a9195c04a47ab4b9.zh-cn.xlf") | Out-Null

$ws2.Range("A2:B5").Font.Underline = $hyperlinkUnderline
$ws2.Range("A2:B5").Font.Color = $hyperlinkColor
$ws2.Range("D2:D5").Font.Underline = $hyperlinkUnderline
$ws2.Range("D2:D5").Font.Color = $hyperlinkColor
$ws2.Range("F2:G2").Font.Underline = $hyperlinkUnderline
$ws2.Range("F2:G2").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(3).Insert()
$ws3.Rows.Item(3).Insert()

$ws3.Range("A3").Value = "42f4766d-43ab-4701-a29c-f352d6e86e17.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "42f4766d-43ab-4701-a29c-f352d6e86e17.83e06547b6f1b7686a5c531cdcbf212a906eda13.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-13 10:41:01"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"

$ws3.Range("A4").Value = "b5df1cc1-dfc1-4509-9a88-1b812e811937.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "b5df1cc1-dfc1-4509-9a88-1b812e811937.e97104f063c9fdfd4fc8a5077fc084d94d44698e.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-13 10:41:01"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("I4").Value = "Include"

$ws3.UsedRange.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b7338fff35708235e71d882e94de97c35b7ecac/e2e/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md", "", "", "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/2b7338fff35708235e71d882e94de97c35b7ecac/e2e/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c0e2050953e0799cc992ba2a50e9ddd1ffea27f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.de-de.xlf", "", "", "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/098ef5c32767506b8a4b577c068dd3d53f16819e/e2e/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md", "", "", "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/af60eb49f1f3f8ac2041bd3491945f825f28db75/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.de-de.xlf", "", "", "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/42f4766d43ab4701a29cf352d6e86e170000000/e2e/42f4766d-43ab-4701-a29c-f352d6e86e17.md", "", "", "42f4766d-43ab-4701-a29c-f352d6e86e17.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/42f4766d43ab4701a29cf352d6e86e170000000/e2e/42f4766d-43ab-4701-a29c-f352d6e86e17.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83e06547b6f1b7686a5c531cdcbf212a906eda13/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/42f4766d-43ab-4701-a29c-f352d6e86e17.83e06547b6f1b7686a5c531cdcbf212a906eda13.de-de.xlf", "", "", "42f4766d-43ab-4701-a29c-f352d6e86e17.83e06547b6f1b7686a5c531cdcbf212a906eda13.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b5df1cc1dfc145099a881b812e8119370000000/e2e/b5df1cc1-dfc1-4509-9a88-1b812e811937.md", "", "", "b5df1cc1-dfc1-4509-9a88-1b812e811937.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/b5df1cc1dfc145099a881b812e8119370000000/e2e/b5df1cc1-dfc1-4509-9a88-1b812e811937.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e97104f063c9fdfd4fc8a5077fc084d94d44698e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b5df1cc1-dfc1-4509-9a88-1b812e811937.e97104f063c9fdfd4fc8a5077fc084d94d44698e.de-de.xlf", "", "", "b5df1cc1-dfc1-4509-9a88-1b812e811937.e97104f063c9fdfd4fc8a5077fc084d94d44698e.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/8b134427431c13d284d55392df0aa999e8a1ea98/e2e/f0ea69a4-53a8-458c-bc21-83bb0ca46deb.md", "", "", "f0ea69a4-53a8-458c-bc21-83bb0ca46deb.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/8b134427431c13d284d55392df0aa999e8a1ea98/e2e/f0ea69a4-53a8-458c-bc21-83bb0ca46deb.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5041ef7b9ad1e13a74db81114bc9b3ee88781c66/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f0ea69a4-53a8-458c-bc21-83bb0ca46deb.91b5d2b6863b784fabb2a985a9195c04a47ab4b9.de-de.xlf", "", "", "f0ea69a4-53a8-458c-bc21-83bb0ca46deb.91b5d2b6863b784fabb2a985a9195c04a47ab4b9.de-de.xlf") | Out-Null

$ws3.Range("A2:B5").Font.Underline = $hyperlinkUnderline
$ws3.Range("A2:B5").Font.Color = $hyperlinkColor
$ws3.Range("D2:D5").Font.Underline = $hyperlinkUnderline
$ws3.Range("D2:D5").Font.Color = $hyperlinkColor
$ws3.Range("F2:G2").Font.Underline = $hyperlinkUnderline
$ws3.Range("F2:G2").Font.Color = $hyperlinkColor

"Report regenerated: added 42f4766d-43ab-4701-a29c-f352d6e86e17 and b5df1cc1-dfc1-4509-9a88-1b812e811937 (Ready for handoff) ahead of f0ea69a4-53a8-458c-bc21-83bb0ca46deb on all sheets."
